$d = $word.ActiveDocument

# Paragraph 5 is "How does the source of positive flows ... Sub-Saharan Africa?"
# (it also carries the trailing " " run + the _GoBack bookmark).
# Paragraph 6 is the trailing empty paragraph at the very end of the body.
# Replace that whole span (paragraph 5 through the end of paragraph 6) with the
# fully-specified target XML in one shot, so the new text / new paragraphs /
# repositioned bookmark / paragraph indent all land exactly as intended.

$startPara = $d.Paragraphs(5)
$endPara = $d.Paragraphs(6)
$target = $d.Range($startPara.Range.Start, $endPara.Range.End)

$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">How does the source of positive flows of agricultural development aid designated for nutrition influence the level of food security in countries of Sub-Saharan Africa? </w:t></w:r></w:p><w:p/><w:p><w:r><w:t>--</w:t></w:r></w:p><w:p><w:r><w:t>How does the source of net flows of agricultural capital flows influence the availability of food in countries of Sub-Saharan Africa?</w:t></w:r></w:p><w:p><w:r><w:t>-Development flows for food security programs</w:t></w:r></w:p><w:p><w:r><w:t>--Bilateral</w:t></w:r></w:p><w:p><w:r><w:t>--Multilateral</w:t></w:r></w:p><w:p><w:r><w:t>--Private</w:t></w:r></w:p><w:p><w:r><w:t>-Agricultural FDI (private)</w:t></w:r></w:p><w:p><w:r><w:t>-Want: remittances</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="720" w:hanging="720"/></w:pPr><w:r><w:t>-Want: trade</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@

$target.InsertXML($xml)
